# Modultafel_Template.xlsx - add new "Titel"/Info/Warning text fields to the
# Einstellungen sheet as a small lookup table (Tabelle4), used by the site
# generator for the page header / footnotes / WPM warning banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Einstellungen")

# --- new header row (G1:K1) ---------------------------------------------
$ws.Range("G1").Value = "Titel"
$ws.Range("H1").Value = "InfoTextOben"
$ws.Range("I1").Value = "WarningTextOben"
$ws.Range("J1").Value = "InfoTextUnten"
$ws.Range("K1").Value = "WarningTextWPM"

# --- new data row (G2:K2) -------------------------------------------------
$ws.Range("G2").Value = "Modultafel Bachelorstudiengang Wirtschaftsinformatik" + [char]10 + "Vollzeit, ab Herbstsemester 2014"
$ws.Range("H2").Value = "Dies ist ein Informations-text Beispiel"
$ws.Range("I2").Value = "Dies ist ein Warnungs-text Beispiel"
$ws.Range("J2").Value = "* = Dieses Modul wird in englischer Sprache durchgeführt" + [char]10 + "Drucktipps: Stellen Sie in der Druckansicht sicher, dass die Option ""Querformat"" aktiviert ist und die Seitenränder in den Seiteneigenschaften auf 0 gesetzt sind."
$ws.Range("K2").Value = "Die Modulbeschreibungen sind aktuell in Überarbeitung und stehen Ihnen ab Montag, 21. November 2022 zur Verfügung."

# Title / the "* = English module" footnote wrap onto multiple lines
$ws.Range("G2").WrapText = $true
$ws.Range("J2").WrapText = $true

# Row heights grow to fit the wrapped text
$ws.Rows.Item(2).RowHeight = 18.2
$ws.Rows.Item(3).RowHeight = 18.2

# --- widen the new columns so the texts are readable ---------------------
$ws.Columns.Item(7).ColumnWidth = 33.125
$ws.Columns.Item(8).ColumnWidth = 33.125
$ws.Columns.Item(10).ColumnWidth = 33.125
$ws.Columns.Item(11).ColumnWidth = 102.875
$ws.Columns.Item(12).ColumnWidth = 33.125

# --- turn G1:K2 into its own table, like the other settings tables -------
$tbl = $ws.ListObjects.Add(1, $ws.Range("G1:K2"), 0, 1)
$tbl.Name = "Tabelle4"
$tbl.TableStyle = "TableStyleLight10"

# --- view tweaks matching the refreshed sheet -----------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws.Range("K20").Select() | Out-Null
